$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.490058565930797
$ws.Range("C2").Value = 1.620881813874908
$ws.Range("B3").Value = 5.110313893513412
$ws.Range("C3").Value = 2.994468826687806
$ws.Range("B4").Value = 6.858508804369644
$ws.Range("C4").Value = 4.759080850396959
$ws.Range("B5").Value = 15.3171969844331
$ws.Range("C5").Value = 6.269858454889558
$ws.Range("B6").Value = 16.33068632295175
$ws.Range("C6").Value = 7.688498984519819
$ws.Range("B7").Value = 20.75638682263203
$ws.Range("C7").Value = 9.070221332874302
$ws.Range("B8").Value = 26.42465011269341
$ws.Range("C8").Value = 10.38584018556291
$ws.Range("B9").Value = 27.09809933776642
$ws.Range("C9").Value = 12.07132950583745
$ws.Range("B10").Value = 35.18829426711937
$ws.Range("C10").Value = 13.31438615145345
$ws.Range("B11").Value = 36.70038112036372
$ws.Range("C11").Value = 14.8015866692244
$ws.Range("B12").Value = 37.30977809531738
$ws.Range("C12").Value = 16.58373071716179
$ws.Range("B13").Value = 38.14229334709679
$ws.Range("C13").Value = 17.82810189444163
$ws.Range("B14").Value = 38.5719663215275
$ws.Range("C14").Value = 19.35076535945387
$ws.Range("B15").Value = 38.84557374366931
$ws.Range("C15").Value = 20.67880509528887
$ws.Range("B16").Value = 39.5243991116746
$ws.Range("C16").Value = 21.97491524829436
$ws.Range("B17").Value = 39.70567846102387
$ws.Range("C17").Value = 23.42679975367606
$ws.Range("B18").Value = 40.00115649419028
$ws.Range("C18").Value = 24.83429714177458
$ws.Range("B19").Value = 42.41216251985806
$ws.Range("C19").Value = 26.56120686319364
$ws.Range("B20").Value = 45.88720228059731
$ws.Range("C20").Value = 27.9409404894594
$ws.Range("B21").Value = 48.40015997679462
$ws.Range("C21").Value = 29.86799222508439
$ws.Range("B22").Value = 48.58488833277983
$ws.Range("C22").Value = 31.04822813559518
$ws.Range("B23").Value = 48.750658882695
$ws.Range("C23").Value = 32.74957969240667
$ws.Range("B24").Value = 50.4278154264795
$ws.Range("C24").Value = 34.16924335174449
$ws.Range("B25").Value = 51.11219023099868
$ws.Range("C25").Value = 35.4819744497531
$ws.Range("B26").Value = 53.43369351216838
$ws.Range("C26").Value = 37.06999772008673
$ws.Range("B27").Value = 54.81225429294069
$ws.Range("C27").Value = 38.30560477136625
$ws.Range("B28").Value = 56.75543399345208
$ws.Range("C28").Value = 39.91400221613734
$ws.Range("B29").Value = 60.02087754636156
$ws.Range("C29").Value = 41.06731740780052
$ws.Range("B30").Value = 60.40734089781912
$ws.Range("C30").Value = 42.77718253795002
$ws.Range("B31").Value = 62.61773297528798
$ws.Range("C31").Value = 44.34640780486198
$ws.Range("B32").Value = 65.82621732973547
$ws.Range("C32").Value = 45.55328779730963
$ws.Range("B33").Value = 67.25424146288067
$ws.Range("C33").Value = 47.30943270671355
$ws.Range("B34").Value = 71.4134597558034
$ws.Range("C34").Value = 48.98347337385841
$ws.Range("B35").Value = 73.05055158660305
$ws.Range("C35").Value = 50.53611865247397
$ws.Range("B36").Value = 73.87543635566935
$ws.Range("C36").Value = 51.74290934553003
$ws.Range("B37").Value = 74.53651078003639
$ws.Range("C37").Value = 53.77212028679437
$ws.Range("B38").Value = 74.67402017275074
$ws.Range("C38").Value = 55.13080165189082
$ws.Range("B39").Value = 79.3898450796483
$ws.Range("C39").Value = 56.67222922117433
$ws.Range("B40").Value = 81.93881362810673
$ws.Range("C40").Value = 58.05102856364998
$ws.Range("B41").Value = 84.00226431520605
$ws.Range("C41").Value = 59.72046871575616
$ws.Range("B42").Value = 84.1231761044728
$ws.Range("C42").Value = 60.845873960382
$ws.Range("B43").Value = 86.09527574473648
$ws.Range("C43").Value = 62.84979759154067
$ws.Range("B44").Value = 86.97679449467499
$ws.Range("C44").Value = 64.31098750533528
$ws.Range("B45").Value = 90.74410627854213
$ws.Range("C45").Value = 66.06866965963884
$ws.Range("B46").Value = 96.34398863897857
$ws.Range("C46").Value = 67.79402899184133
$ws.Range("B47").Value = 98.46719278375217
$ws.Range("C47").Value = 69.24081036578599

# Remove now-unused rows 48 and 49
$ws.Range("A48:C49").EntireRow.Delete()

